# Insert a new weekly record at row 64 (shifting the existing rows 64-98
# down to 65-99) for "Terminal La Palmera de La Serena - Jengibre".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(64).Insert()

$ws.Range("A64").Value2 = 8
$ws.Range("B64").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C64").Value2 = "Coquimbo"
$ws.Range("D64").Value2 = 44942
$ws.Range("E64").Value2 = 4
$ws.Range("F64").Value2 = 100114007
$ws.Range("G64").Value2 = "Jengibre"
$ws.Range("H64").Value2 = "Sin especificar"
$ws.Range("I64").Value2 = "Primera"
$ws.Range("J64").Value2 = 440
$ws.Range("K64").Value2 = 14000
$ws.Range("L64").Value2 = 15000
$ws.Range("M64").Value2 = 14500
$ws.Range("N64").Value2 = "$/caja 13 kilos"
$ws.Range("O64").Value2 = "Perú"
$ws.Range("P64").Value2 = 1115
$ws.Range("Q64").Value2 = 13
$ws.Range("R64").Value2 = "Hortaliza"
